# Update "想去人数" (number of people interested) figures for a handful of
# events that appear on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 11697
$wsExhibit.Range("F7").Value  = 11652
$wsExhibit.Range("F8").Value  = 482
$wsExhibit.Range("F12").Value = 5767
$wsExhibit.Range("F14").Value = 3516

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 11697
$wsAll.Range("F9").Value  = 11652
$wsAll.Range("F10").Value = 482
$wsAll.Range("F15").Value = 5767
$wsAll.Range("F17").Value = 3516
